# Update cryptocurrency price/volume figures per the latest GitHub Actions scrape.
# D = Price column (free-text, "." used as thousands+decimal separators in source feed)
# E = Volume(1h) column (free-text percentage, padded with two leading/trailing spaces)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.921.27"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "2.585.10"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'582.66"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'147.00"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E8").Value = "  +2.01%  "
$ws.Range("E9").Value = "  +2.32%  "
$ws.Range("D10").Value = "'5.66"
$ws.Range("E10").Value = "  +2.39%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "'27.40"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").Value = "3.047.92"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("D15").Value = "62.800.50"
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("E16").Value = "  +2.79%  "
$ws.Range("D17").Value = "2.590.84"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").Value = "'342.07"
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "'67.04"
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("D24").Value = "2.708.95"
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("D26").Value = "'1.59"
$ws.Range("E26").Value = "  -2.13%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "'8.33"
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("E29").Value = "  +6.25%  "
$ws.Range("E30").Value = "  -2.15%  "
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").Value = "'469.90"
$ws.Range("E32").Value = "  +13.91%  "
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("D34").Value = "'176.08"
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("E35").Value = "  +3.68%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("D38").Value = "'19.05"
$ws.Range("E38").Value = "  -0.57%  "
$ws.Range("E39").Value = "  +3.89%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("D42").Value = "'157.53"
$ws.Range("E42").Value = "  +4.10%  "
$ws.Range("D43").Value = "'3.75"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("E44").Value = "  +4.94%  "
$ws.Range("D45").Value = "'21.19"
$ws.Range("E45").Value = "  +1.26%  "
$ws.Range("D46").Value = "'0.0541"
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("D47").Value = "'0.0966"
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("E48").Value = "  -0.96%  "
$ws.Range("D49").Value = "'18.35"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("D51").Value = "'11.42"
$ws.Range("E51").Value = "  +1.07%  "

# The prices above that are plain decimal numbers (e.g. "147.00") get auto-parsed
# by Excel as numeric input when quote-prefixed for text entry, which also stamps the
# cell with a quote-prefix style. Reset those cells back to the default "Normal" style
# so formatting matches the untouched cells exactly (only the text content changed).
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
